# Auto-generated Excel COM-interop script applying the scheduled market-data
# refresh to each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2593.275
$ws.Range("I132").Value = 1816.28
$ws.Range("J132").Value = 3888.2666
$ws.Range("K132").Value = 5448.84
$ws.Range("L132").Value = 11664.7998
$ws.Range("M132").Value = -2918.84
$ws.Range("N132").Value = -16724.7998

$ws.Range("H138").Value = 3080091.2
$ws.Range("I138").Value = 11112782
$ws.Range("J138").Value = 3741.7234
$ws.Range("K138").Value = 33338346
$ws.Range("L138").Value = 11225.1702
$ws.Range("M138").Value = -33333206
$ws.Range("N138").Value = -21505.1702

$ws.Range("H139").Value = 61238.668
$ws.Range("J139").Value = 61238.668
$ws.Range("L139").Value = 61238.668
$ws.Range("N139").Value = -71518.66800000001


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1488.375
$ws.Range("I61").Value = 1146.3636
$ws.Range("J61").Value = 2240.8
$ws.Range("K61").Value = 1146.3636
$ws.Range("L61").Value = 2240.8
$ws.Range("M61").Value = -934.3635999999999
$ws.Range("N61").Value = -2664.8

$ws.Range("H74").Value = 1199.1428
$ws.Range("I74").Value = 1002.6667
$ws.Range("J74").Value = 1552.8
$ws.Range("K74").Value = 1002.6667
$ws.Range("L74").Value = 1552.8
$ws.Range("M74").Value = -128.6667
$ws.Range("N74").Value = -3300.8

$ws.Range("H77").Value = 1199.1428
$ws.Range("I77").Value = 1002.6667
$ws.Range("J77").Value = 1552.8
$ws.Range("K77").Value = 5013.3335
$ws.Range("L77").Value = 7764
$ws.Range("M77").Value = -645.3334999999997
$ws.Range("N77").Value = -16500

$ws.Range("H97").Value = 360.29413
$ws.Range("I97").Value = 364.0625
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 364.0625
$ws.Range("L97").Value = 300
$ws.Range("M97").Value = 131.9375
$ws.Range("N97").Value = -1292

$ws.Range("H123").Value = 40429
$ws.Range("J123").Value = 40429
$ws.Range("L123").Value = 40429
$ws.Range("N123").Value = -50229

$ws.Range("H136").Value = 1488.375
$ws.Range("I136").Value = 1146.3636
$ws.Range("J136").Value = 2240.8
$ws.Range("K136").Value = 3439.0908
$ws.Range("L136").Value = 6722.400000000001
$ws.Range("M136").Value = -889.0907999999999
$ws.Range("N136").Value = -11822.4


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 74739.5
$ws.Range("I57").Value = 50709
$ws.Range("K57").Value = 50709
$ws.Range("M57").Value = -49989

$ws.Range("H94").Value = 1509
$ws.Range("I94").Value = 1151
$ws.Range("J94").Value = 2225
$ws.Range("K94").Value = 1151
$ws.Range("L94").Value = 2225
$ws.Range("M94").Value = -700
$ws.Range("N94").Value = -3127

$ws.Range("H99").Value = 1581.24
$ws.Range("I99").Value = 1372.8572
$ws.Range("J99").Value = 1846.4546
$ws.Range("K99").Value = 1372.8572
$ws.Range("L99").Value = 1846.4546
$ws.Range("M99").Value = 125.1428000000001
$ws.Range("N99").Value = -4842.4546

$ws.Range("H134").Value = 2071.6206
$ws.Range("I134").Value = 1889.3
$ws.Range("J134").Value = 2476.7778
$ws.Range("K134").Value = 5667.9
$ws.Range("L134").Value = 7430.3334
$ws.Range("M134").Value = -3132.9
$ws.Range("N134").Value = -12500.3334

$ws.Range("H136").Value = 74739.5
$ws.Range("I136").Value = 50709
$ws.Range("K136").Value = 50709
$ws.Range("M136").Value = -45609


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1192.2413
$ws.Range("I58").Value = 1191.0416
$ws.Range("J58").Value = 1198
$ws.Range("K58").Value = 1191.0416
$ws.Range("L58").Value = 1198
$ws.Range("M58").Value = -988.0416
$ws.Range("N58").Value = -1604

$ws.Range("H132").Value = 2207.6924
$ws.Range("I132").Value = 1850.25
$ws.Range("J132").Value = 2779.6
$ws.Range("K132").Value = 5550.75
$ws.Range("L132").Value = 8338.799999999999
$ws.Range("M132").Value = -3020.75
$ws.Range("N132").Value = -13398.8

$ws.Range("H134").Value = 1375.3823
$ws.Range("I134").Value = 1239.2069
$ws.Range("J134").Value = 2165.2
$ws.Range("K134").Value = 3717.620699999999
$ws.Range("L134").Value = 6495.599999999999
$ws.Range("M134").Value = -1182.620699999999
$ws.Range("N134").Value = -11565.6

$ws.Range("H136").Value = 1192.2413
$ws.Range("I136").Value = 1191.0416
$ws.Range("J136").Value = 1198
$ws.Range("K136").Value = 3573.1248
$ws.Range("L136").Value = 3594
$ws.Range("M136").Value = -1023.1248
$ws.Range("N136").Value = -8694


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 3042
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 3465.6667
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 10397.0001
$ws.Range("M10").Value = -1361
$ws.Range("N10").Value = -10675.0001

$ws.Range("H113").Value = 672.6667
$ws.Range("I113").Value = 640
$ws.Range("J113").Value = 689
$ws.Range("K113").Value = 1920
$ws.Range("L113").Value = 2067
$ws.Range("M113").Value = 250
$ws.Range("N113").Value = -6407

$ws.Range("H139").Value = 1915.174
$ws.Range("I139").Value = 1225
$ws.Range("J139").Value = 2988.7778
$ws.Range("K139").Value = 3675
$ws.Range("L139").Value = 8966.3334
$ws.Range("M139").Value = 1465
$ws.Range("N139").Value = -19246.3334

$ws.Range("H140").Value = 2463.8262
$ws.Range("I140").Value = 861.25
$ws.Range("J140").Value = 4212.091
$ws.Range("K140").Value = 2583.75
$ws.Range("L140").Value = 12636.273
$ws.Range("M140").Value = 2596.25
$ws.Range("N140").Value = -22996.273

$ws.Range("H141").Value = 4710.1
$ws.Range("I141").Value = 3350
$ws.Range("K141").Value = 10050
$ws.Range("M141").Value = -4870


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4727330.5
$ws.Range("I14").Value = 5901663
$ws.Range("J14").Value = 30000
$ws.Range("K14").Value = 5901663
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = -5901495
$ws.Range("N14").Value = -30336

$ws.Range("H22").Value = 3000
$ws.Range("J22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -4058

$ws.Range("H57").Value = 39949.25
$ws.Range("J57").Value = 39949.25
$ws.Range("L57").Value = 39949.25
$ws.Range("N57").Value = -41589.25

$ws.Range("H97").Value = 19488.555
$ws.Range("I97").Value = 22655.652
$ws.Range("J97").Value = 1277.75
$ws.Range("K97").Value = 22655.652
$ws.Range("L97").Value = 1277.75
$ws.Range("M97").Value = -22159.652
$ws.Range("N97").Value = -2269.75

$ws.Range("H122").Value = 3162.1428
$ws.Range("I122").Value = 3024.8125
$ws.Range("J122").Value = 3601.6
$ws.Range("K122").Value = 9074.4375
$ws.Range("L122").Value = 10804.8
$ws.Range("M122").Value = -6624.4375
$ws.Range("N122").Value = -15704.8

$ws.Range("H123").Value = 30325.2
$ws.Range("J123").Value = 30325.2
$ws.Range("L123").Value = 30325.2
$ws.Range("N123").Value = -35225.2

$ws.Range("H126").Value = 4205.4
$ws.Range("I126").Value = 3999
$ws.Range("K126").Value = 11997
$ws.Range("M126").Value = -9527

$ws.Range("H132").Value = 2329.95
$ws.Range("I132").Value = 1478.7142
$ws.Range("J132").Value = 4316.1665
$ws.Range("K132").Value = 4436.142599999999
$ws.Range("L132").Value = 12948.4995
$ws.Range("M132").Value = -1906.142599999999
$ws.Range("N132").Value = -18008.4995


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3266.4443
$ws.Range("I40").Value = 2860.111
$ws.Range("J40").Value = 3672.7778
$ws.Range("K40").Value = 2860.111
$ws.Range("L40").Value = 3672.7778
$ws.Range("M40").Value = -2724.111
$ws.Range("N40").Value = -3944.7778

$ws.Range("H46").Value = 1270
$ws.Range("I46").Value = 1120
$ws.Range("K46").Value = 1120
$ws.Range("M46").Value = -932

$ws.Range("H96").Value = 30193
$ws.Range("J96").Value = 30193
$ws.Range("L96").Value = 30193
$ws.Range("N96").Value = -35685

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 7000
$ws.Range("J19").Value = 7000
$ws.Range("L19").Value = 7000
$ws.Range("N19").Value = -7348

$ws.Range("H54").Value = 26300
$ws.Range("J54").Value = 26300
$ws.Range("L54").Value = 26300
$ws.Range("N54").Value = -27340

$ws.Range("H132").Value = 2910.7
$ws.Range("I132").Value = 2020.8
$ws.Range("J132").Value = 3800.6
$ws.Range("K132").Value = 6062.4
$ws.Range("L132").Value = 11401.8
$ws.Range("M132").Value = -3532.4
$ws.Range("N132").Value = -16461.8

